# Brainstorm.docx edit: append a new bullet under the "Ta scheduling" /
# "refinements" list, right after the final existing paragraph
# ("But we could let them explicitly say I don't care as much about the
# previous step"). The new bullet sits one level up from that paragraph
# (same list, ilvl 2 / "List Paragraph" style, numId 3 continued).

$d = $word.ActiveDocument

# Grab the last paragraph in the document (the "But we could..." bullet).
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)

# Insert a new paragraph right after it; Word clones the preceding
# paragraph's list/style formatting onto the new one.
$lastPara.Range.InsertParagraphAfter()

# Re-fetch the freshly created paragraph by index (the previous object
# references are not guaranteed to reflect the mutated document) and
# give it the text + promote it from ilvl 3 up to ilvl 2 (ListLevelNumber
# is 1-based in the object model, so level 3 == w:ilvl 2).
$newIndex = $d.Paragraphs.Count
$newPara = $d.Paragraphs($newIndex)
$newPara.Range.ListFormat.ListLevelNumber = 3
$newPara.Range.Text = "One way to make the user do this is by building constraints like nodes in a tree, and the user can have different branches that try different things"
